# Update countries & provincias Spain
#
# Refreshes the COVID-19 snapshot on the "Pais" sheet to the 13:22 update
# (previous snapshot was 12:52). A handful of countries' case counts grew
# enough to change their ranking in the (descending, by "Casos totales")
# list, so besides the raw numbers, a few rows' country names (column A)
# shift down to the next-lower rank and the displaced countries - Bielorrusia
# and Senegal - move up into the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header timestamp
$ws.Range("A1").Value = "Datos actualizados a 10 de Abril de 2020 a las 13:22"

# Iran (row 10)
$ws.Range("B10").Value = 68192
$ws.Range("C10").Value = 1972
$ws.Range("D10").Value = 35465
$ws.Range("E10").Value = 28495
$ws.Range("F10").Value = 3969
$ws.Range("G10").Value = 122
$ws.Range("H10").Value = 4232

# Austria (row 19)
$ws.Range("B19").Value = 13431
$ws.Range("C19").Value = 187
$ws.Range("E19").Value = 7048

# Rows 52-59: Bielorrusia's case count overtakes Grecia/Sudafrica/Singapur/
# Argentina/Egipto/Argelia/Islandia, so each row's country label shifts down
# one rank and gets the data for its new occupant.
$ws.Range("A52").Value = "Bielorrusia"
$ws.Range("B52").Value = 1981
$ws.Range("C52").Value = 495
$ws.Range("D52").Value = 169
$ws.Range("E52").Value = 1793
$ws.Range("F52").Value = 72
$ws.Range("G52").Value = 3
$ws.Range("H52").Value = 19

$ws.Range("A53").Value = "Grecia"
$ws.Range("B53").Value = 1955
$ws.Range("D53").Value = 269
$ws.Range("E53").Value = 1599
$ws.Range("F53").Value = 79
$ws.Range("H53").Value = 87

$ws.Range("A54").Value = "Sudafrica"
$ws.Range("B54").Value = 1934
$ws.Range("D54").Value = 95
$ws.Range("E54").Value = 1821
$ws.Range("F54").Value = 7
$ws.Range("H54").Value = 18

$ws.Range("A55").Value = "Singapur"
$ws.Range("B55").Value = 1910
$ws.Range("D55").Value = 460
$ws.Range("E55").Value = 1444
$ws.Range("F55").Value = 29
$ws.Range("H55").Value = 6

$ws.Range("A56").Value = "Argentina"
$ws.Range("B56").Value = 1894
$ws.Range("D56").Value = 365
$ws.Range("E56").Value = 1450
$ws.Range("F56").Value = 96
$ws.Range("H56").Value = 79

$ws.Range("A57").Value = "Egipto"
$ws.Range("B57").Value = 1699
$ws.Range("D57").Value = 348
$ws.Range("E57").Value = 1233
$ws.Range("F57").Value = 0
$ws.Range("H57").Value = 118

$ws.Range("A58").Value = "Argelia"
$ws.Range("B58").Value = 1666
$ws.Range("D58").Value = 347
$ws.Range("E58").Value = 1084
$ws.Range("F58").Value = 46
$ws.Range("H58").Value = 235

$ws.Range("A59").Value = "Islandia"
$ws.Range("B59").Value = 1648
$ws.Range("D59").Value = 688
$ws.Range("E59").Value = 954
$ws.Range("F59").Value = 11
$ws.Range("H59").Value = 6

# Uzbekistan (row 81)
$ws.Range("D81").Value = 42
$ws.Range("E81").Value = 579

# Mauricio (row 104)
$ws.Range("E104").Value = 282
$ws.Range("G104").Value = 2
$ws.Range("H104").Value = 9

# Estado de Palestina (row 108)
$ws.Range("D108").Value = 45
$ws.Range("E108").Value = 219
$ws.Range("G108").Value = 1
$ws.Range("H108").Value = 2

# Rows 109-111: Senegal's case count overtakes Montenegro/Vietnam, so those
# two ranks shift down and Senegal moves up into row 109.
$ws.Range("A109").Value = "Senegal"
$ws.Range("B109").Value = 265
$ws.Range("C109").Value = 15
$ws.Range("D109").Value = 137
$ws.Range("E109").Value = 126
$ws.Range("F109").Value = 1

$ws.Range("A110").Value = "Montenegro"
$ws.Range("C110").Value = 3
$ws.Range("D110").Value = 4
$ws.Range("E110").Value = 249
$ws.Range("F110").Value = 7
$ws.Range("H110").Value = 2

$ws.Range("A111").Value = "Vietnam"
$ws.Range("B111").Value = 255
$ws.Range("D111").Value = 144
$ws.Range("E111").Value = 111
$ws.Range("F111").Value = 8
$ws.Range("H111").Value = 0
